# Rename the embedded-picture display names:
#   - Header logo  (BTec_Logo-Orange):        image1.jpg -> image2.jpg
#   - Footer logos (Pearson PowerPoint logo):  image2.png -> image1.png  (x2, one per footer)
#
# These are the wp:docPr/pic:cNvPr "name" attributes on the inline pictures
# living in the document's header/footer parts - not document body text, so
# Find/Replace can't reach them. We walk Sections -> Headers/Footers ->
# Range.InlineShapes and rename each picture via the InlineShape's Name
# property, matching on its (unique, unchanged) AlternativeText/descr so the
# right picture is renamed even if shape ordering ever shifts.

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    for ($hi = 1; $hi -le 3; $hi++) {
        $hdr = $sec.Headers.Item($hi)
        $shapes = $hdr.Range.InlineShapes
        for ($k = 1; $k -le $shapes.Count; $k++) {
            $shp = $shapes.Item($k)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange" -and $shp.Name -eq "image1.jpg") {
                $shp.Name = "image2.jpg"
            }
        }
    }

    for ($fi = 1; $fi -le 3; $fi++) {
        $ftr = $sec.Footers.Item($fi)
        $shapes = $ftr.Range.InlineShapes
        for ($k = 1; $k -le $shapes.Count; $k++) {
            $shp = $shapes.Item($k)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" -and $shp.Name -eq "image2.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}
